$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted cells stay as text (they store numbers/percentages as strings)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "312.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "8.41%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "9.19%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.347"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "4.93%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07577"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "13.37%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.844"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "7.13%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.708"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "8.81%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.606"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "17.33%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9132"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.50%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01688"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2,504.56%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1721"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "8.36%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07771"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "15.24%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08241"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "8.93%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03025"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "3.02%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09875"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "9.99%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001525"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.69%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04551"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.91%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006518"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.70%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.501"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.62%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.242"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.90%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3311"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.07%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1317"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.65%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.174"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.40%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1620"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.39%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001217"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.14%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004498"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.38%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001297"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "8.14%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001739"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "7.54%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04606"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "8.29%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007203"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "7.25%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "10.45%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002255"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.27%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01403"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.43%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006150"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.68%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.82%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01398"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "7.06%"
